# edit.ps1
# Applies the commit "Add files via upload":
#  - Renames "Лист1" -> "ссылки" (sheet1) and fills it with a genre -> Apple Music
#    playlist link lookup table (25 rows, columns A:B), makes it the active sheet.
#  - Leaves "стикеры" (sheet2) content as-is but updates its view/selection and
#    resets column C's width back to the (default) width.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename first sheet ---
$ws1.Name = "ссылки"

# --- Populate "ссылки" sheet with header + 25 rows of genre/link pairs ---
$rows = @(
    @(1, 'ключевое слово (жанр)', 'ссылка'),
    @(2, 'джаз', 'https://music.apple.com/ru/playlist/сонный-джаз/pl.f31fb34cf8ad46f0a27c960fa805cc4c'),
    @(3, 'рок', 'https://music.apple.com/ru/playlist/пульс-рока/pl.aea893d092ed40a393d217ced33bde62'),
    @(4, 'инди', 'https://music.apple.com/ru/playlist/под-звездным-небом/pl.cd1e46899e31421285d7e655bac51cbd'),
    @(5, 'хип-хоп', 'https://music.apple.com/ru/playlist/rap-life/pl.abe8ba42278f4ef490e3a9fc5ec8e8c5'),
    @(6, 'классика', 'https://music.apple.com/ru/playlist/музыка-без-слов/pl.60280cd08f2e4e1fa2cfb27a8466f9c5'),
    @(7, 'блюз', 'https://music.apple.com/ru/playlist/блюз-с-большой-дороги/pl.a9faca07cf8f47e19f1819b0f5a2e765'),
    @(8, 'поп', 'https://music.apple.com/ru/playlist/новинки-поп-музыки/pl.5ee8333dbe944d9f9151e97d92d1ead9'),
    @(9, 'метал', 'https://music.apple.com/ru/playlist/экстремальный-метал/pl.1baada6675ca477cbe9946b3d21c5757'),
    @(10, 'hard рок', 'https://music.apple.com/ru/playlist/зловещее-звучание/pl.f72bb902f8684581bb26760069b50e43'),
    @(11, 'ретро', 'https://music.apple.com/ru/playlist/motown-главное/pl.15378bf1b2624a56af9751390242a19d'),
    @(12, 'регги', 'https://music.apple.com/ru/playlist/регги-вайб/pl.e75fb4f0f6f649a89f7c28ef4cc0442f'),
    @(13, 'фолк', 'https://music.apple.com/ru/playlist/фолк-главное/pl.ced4e8788cab46e7982ba4a26e5211a7'),
    @(14, 'панк рок', 'https://music.apple.com/ru/playlist/проснись-и-вой/pl.937d9aefc9da498aaa627e7b62e318b4'),
    @(15, 'альтернатива', 'https://music.apple.com/ru/playlist/пульс-альтернативы/pl.132b9a231cbf464086b4f838b2726f94'),
    @(16, 'электроника', 'https://music.apple.com/ru/playlist/бесконечный-бит/pl.4705ab1ed97c4f4bb54f48940faf5623'),
    @(17, 'вечеринка', 'https://music.apple.com/ru/playlist/коктейльная-вечеринка/pl.d33b5bd820cf47ffb87889db225bd943'),
    @(18, 'хорошее настроение ', 'https://music.apple.com/ru/playlist/хорошее-настроение/pl.10fc76a3edc14e759deb60535854e339'),
    @(19, 'сон', 'https://music.apple.com/ru/playlist/спокойный-сон/pl.f31a7e6b60ab4e0d995d837db65bbf00'),
    @(20, 'концентрация', 'https://music.apple.com/ru/playlist/полная-концентрация/pl.a4e197979fc74b2a91b3cdf869f12aa5'),
    @(21, 'спорт ', 'https://music.apple.com/ru/playlist/только-тренировка/pl.ad0ee1557e3e4feba314fd70f7982766'),
    @(22, 'спокойствие ', 'https://music.apple.com/ru/playlist/только-умиротворение/pl.ffc344338c3d4ff394ddcf94d766c143'),
    @(23, 'романтика', 'https://music.apple.com/ru/playlist/любовь/pl.7d7525d9145c4aba9a6295753505c98d'),
    @(24, 'меланхолия ', 'https://music.apple.com/ru/playlist/в-изоляции/pl.464c6868d80d4c0891e22568c234829f'),
    @(25, 'мотивация', 'https://music.apple.com/ru/playlist/только-мотивация/pl.047294ae14a24e5993d1f7ab2b127188'),
)

foreach ($row in $rows) {
    $r = $row[0]
    $keyword = $row[1]
    $link = $row[2]
    $ws1.Cells.Item($r, 1).Value = $keyword
    $ws1.Cells.Item($r, 2).Value = $link
}

# --- Page setup for "ссылки" (paperSize=9 / A4, portrait) ---
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- "стикеры" sheet: reset column C width back to the default width ---
$ws2.Columns.Item(3).ColumnWidth = 8.25

# --- Selection / active-sheet bookkeeping ---
# Select on "стикеры" first, then "ссылки" last so "ссылки" ends up as the
# active (tabSelected) sheet / tab, matching the target workbook view.
$ws2.Range("G5").Select()
$ws1.Range("A16:XFD16").Select()

Write-Host "edit complete"
